# Pill3 for Beriozka prepared
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the input (non-formula) cells; dependent formula cells
# (C6, D6, C7, D7, D8, B9, ...) recalculate automatically.
$ws.Range("B4").Value = 3
$ws.Range("B6").Value = 901
$ws.Range("B7").Value = 0

$excel.Calculate()
